# Updated symbol list - refresh crypto price snapshot values (Price column)
# and a couple of "Worstin24h" label tweaks on sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores its values as plain text in this workbook
# (not as numbers), so we prefix with a single quote to force Excel to keep
# them as text instead of auto-converting to numeric cells.
$ws.Range("D2").Value  = "'245.67"
$ws.Range("D3").Value  = "'24.02"
$ws.Range("D4").Value  = "'5.209"
$ws.Range("D5").Value  = "'0.05786"
$ws.Range("D6").Value  = "'6.515"
$ws.Range("D7").Value  = "'3.122"
$ws.Range("D8").Value  = "'0.8158"
$ws.Range("D9").Value  = "'0.8507"
$ws.Range("D10").Value = "'0.1356"
$ws.Range("D11").Value = "'0.06953"
$ws.Range("D12").Value = "'0.03152"
$ws.Range("D13").Value = "'0.02880"
$ws.Range("D14").Value = "'0.09384"
$ws.Range("D15").Value = "'3.761"
$ws.Range("D16").Value = "'0.001512"
$ws.Range("D17").Value = "'0.04676"
$ws.Range("D18").Value = "'0.0006003"
$ws.Range("D19").Value = "'0.006274"
$ws.Range("D21").Value = "'0.004301"
$ws.Range("D22").Value = "'0.00008602"
$ws.Range("D23").Value = "'3.499"
$ws.Range("D25").Value = "'0.3173"
$ws.Range("D28").Value = "'0.0002330"
$ws.Range("D40").Value = "'0.03649"
$ws.Range("D41").Value = "'0.006281"
$ws.Range("D42").Value = "'0.1052"
$ws.Range("D43").Value = "'0.002851"
$ws.Range("D44").Value = "'0.007486"
$ws.Range("D45").Value = "'0.00005282"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D47").Value = "'0.2802"
$ws.Range("D48").Value = "'0.002339"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D50").Value = "'0.0002001"

# Row 47/48 label (E column) also changed: the "Worstin24h" marker moved
# from the BOLO row to the CoinbaseStockToken row.
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("E48").Value = "47BOLOBOLO"

Write-Host "Applied crypto price/label refresh to sheet1"
